$wb = $excel.ActiveWorkbook

# Insert a new "is_targeted list" sheet right before "ms_source list"
# (so the sheet order becomes: ... analyte_class list, is_targeted list,
#  ms_source list, polarity list, resolution_x_unit list, resolution_y_unit list)
$msSource = $wb.Worksheets.Item("ms_source list")
$isTargeted = $wb.Worksheets.Add($msSource)
$isTargeted.Name = "is_targeted list"

# Populate it with the two boolean-as-text options, forced to text via a
# leading apostrophe so they are stored as shared strings ("TRUE"/"FALSE")
# rather than native boolean cells. Reset the style afterwards so the
# quote-prefix marker doesn't leave a stray cell style behind.
$isTargeted.Range("A1").Value = "'TRUE"
$isTargeted.Range("A2").Value = "'FALSE"
$isTargeted.Range("A1:A2").Style = "Normal"

# Point the N column (is_targeted) validation at the new list sheet instead
# of the inline "TRUE,FALSE" formula, and update the error text/title to
# match the other list-based validations.
$ws1 = $wb.Worksheets.Item("Export as TSV")
$v = $ws1.Range("N2:N1048576").Validation
$v.Modify(3, 1, 1, "='is_targeted list'!`$A`$1:`$A`$2")
$v.ErrorTitle = "Value must come from list"
$v.ErrorMessage = "Value must be one of: TRUE / FALSE."

# Restore the originally active sheet so we don't leave the new helper
# sheet selected as a side effect of having just created it.
$ws1.Activate()
